$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genres = @(
    "Drum & Bass",
    "House",
    "Trance",
    "Techno",
    "Nu Disco/Indie Dance",
    "Melodic Progressive",
    "Twerk",
    "Ambient",
    "Grime",
    "BreakBeat/Breaks",
    "Progressive House",
    "Wave",
    "Progressive Trance",
    "G-House",
    "Tech House",
    "Deep House"
)

for ($i = 0; $i -lt $genres.Length; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = $genres[$i]
}

$ws.Range("H9").Select()
